# Apply "Misc card changes" edit: Challenge card rework.
# Rewords trigger timing from "交锋时" (on clash) to "开战时" (on battle start)
# for two effect texts, and turns a third effect into an "end of round" trigger.
# Also moves the sheet selection from D15 to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (落穴 / Pitfall): effect text timing changed from "交锋时" to "回合结束时"
$ws.Range("D4").Value = "回合结束时：如果本牌所在槽位和对位槽位的怪物牌合计数量不小于本牌点数，则将那些怪物牌全部消灭，然后消灭本牌。"

# Row 3 (流沙 / Quicksand): effect text timing changed from "交锋时" to "开战时"
$ws.Range("D3").Value = "开战时：本牌所在槽位和对位槽位的所有怪物牌点数变为1。"

# Row 9 (传送阵 / Teleporter): effect text timing changed from "交锋时" to "开战时"
$ws.Range("D9").Value = "开战时：如果本牌所在槽位和对位槽位的怪物牌合计数量大于1，则将那些怪物牌洗回主牌堆，然后消灭本牌。"

# Update the active selection to match the saved state in the new file (D10).
$ws.Range("D10").Select()
